$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Australia ALeague")

# --- Row 104 / 105 swap of match data (ids 7127370 / 7127374) ---
# Row 104 becomes what used to be row 105's data
$ws.Range("B104").Value = 7127374
$ws.Range("F104").Value = "Central Coast Mariners"
$ws.Range("G104").Value = "Western Sydney Wanderers"
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = "H"
$ws.Range("K104").Value = 1.909
$ws.Range("M104").Value = 3.6
$ws.Range("N104").Value = 2.15
$ws.Range("O104").Value = 3.6
$ws.Range("P104").Value = 3.25
$ws.Range("Q104").Value = -0.25
$ws.Range("R104").Value = 1.86
$ws.Range("S104").Value = 2.04
$ws.Range("T104").Value = 2.75
$ws.Range("U104").Value = 1.975
$ws.Range("V104").Value = 1.875
$ws.Range("W104").Value = 1.15
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 0.8600000000000001
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = -1
$ws.Range("AC104").Value = 0.875

# Row 105 becomes what used to be row 104's data
$ws.Range("B105").Value = 7127370
$ws.Range("F105").Value = "Macarthur FC"
$ws.Range("G105").Value = "Wellington Phoenix"
$ws.Range("I105").Value = 2
$ws.Range("J105").Value = "A"
$ws.Range("K105").Value = 2.4
$ws.Range("M105").Value = 2.625
$ws.Range("N105").Value = 2.375
$ws.Range("O105").Value = 3.8
$ws.Range("P105").Value = 2.75
$ws.Range("Q105").Value = 0
$ws.Range("R105").Value = 1.8
$ws.Range("S105").Value = 2.05
$ws.Range("T105").Value = 3
$ws.Range("U105").Value = 1.9
$ws.Range("V105").Value = 1.95
$ws.Range("W105").Value = -1
$ws.Range("Y105").Value = 1.75
$ws.Range("Z105").Value = -1
$ws.Range("AA105").Value = 1.05
$ws.Range("AB105").Value = 0
$ws.Range("AC105").Value = -0

# --- Row 138: add FTHG/FTAG/FTR and update odds columns ---
$ws.Range("H138").Value = 2
$ws.Range("I138").Value = 1
$ws.Range("J138").Value = "H"
$ws.Range("N138").Value = 2.4
$ws.Range("O138").Value = 3.6
$ws.Range("P138").Value = 2.8
$ws.Range("Q138").Value = 0
$ws.Range("R138").Value = 1.83
$ws.Range("S138").Value = 2.07
$ws.Range("T138").Value = 2.75
$ws.Range("U138").Value = 2.025
$ws.Range("V138").Value = 1.825
$ws.Range("W138").Value = 1.4
$ws.Range("X138").Value = -1
$ws.Range("Y138").Value = -1
$ws.Range("Z138").Value = 0.8300000000000001
$ws.Range("AA138").Value = -1
$ws.Range("AB138").Value = 0.5125
$ws.Range("AC138").Value = -0.5

# --- Row 139 ---
$ws.Range("N139").Value = 4
$ws.Range("P139").Value = 1.75
$ws.Range("R139").Value = 1.87
$ws.Range("S139").Value = 2.03
$ws.Range("U139").Value = 1.825
$ws.Range("V139").Value = 2.025

# --- Row 140 ---
$ws.Range("U140").Value = 2
$ws.Range("V140").Value = 1.85

# --- Row 141 ---
$ws.Range("N141").Value = 2.05
$ws.Range("O141").Value = 3.8
$ws.Range("R141").Value = 1.86
$ws.Range("S141").Value = 2.04
$ws.Range("U141").Value = 1.825
$ws.Range("V141").Value = 2.025

# --- Row 142 ---
$ws.Range("N142").Value = 1.666
$ws.Range("O142").Value = 3.75
$ws.Range("P142").Value = 5
$ws.Range("R142").Value = 1.97
$ws.Range("S142").Value = 1.93
$ws.Range("U142").Value = 1.975
$ws.Range("V142").Value = 1.875

# --- Row 143 ---
$ws.Range("N143").Value = 2.2
$ws.Range("R143").Value = 2
$ws.Range("S143").Value = 1.9
$ws.Range("U143").Value = 1.825
$ws.Range("V143").Value = 2.025

# --- Row 144 ---
$ws.Range("R144").Value = 2
$ws.Range("S144").Value = 1.9
